$wb = $excel.ActiveWorkbook
$wsFilepath = $wb.Worksheets.Item("Filepath")

# --- Filepath sheet: rename the generation/exchange data directories ---
$wsFilepath.Range("B2").Value = "../Data_Generation/"
$wsFilepath.Range("B3").Value = "../Data_Exchanges/"

# --- Insert two new rows above the old "mapping file" row (7 & 8), pushing
#     "mapping file" and everything below it down by two rows ---
$wsFilepath.Range("A7:A8").EntireRow.Insert()

# New rows 7 & 8: separate "saving generation" / "saving exchanges" labels
$wsFilepath.Range("A7").Value = "saving generation"
$wsFilepath.Range("A8").Value = "saving exchanges"

# The old "mapping file" row (now row 10) no longer carries a default path
$wsFilepath.Range("B10").Value = ""

# --- Selection / active sheet bookkeeping ---
$wsFilepath.Range("B10").Select()
$wsFilepath.Activate()
